$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 7.281461919790895
$ws.Cells.Item(2, 4).Value = 10.22403199079876
$ws.Cells.Item(2, 5).Value = 13.9548778311738
$ws.Cells.Item(2, 6).Value = 27.25277034875602
$ws.Cells.Item(2, 7).Value = 26.51541891655085
$ws.Cells.Item(2, 8).Value = 13.1576907187012
$ws.Cells.Item(2, 9).Value = 17.82547295445669
$ws.Cells.Item(2, 10).Value = 9.79760118383477
$ws.Cells.Item(2, 13).Value = 17.93939382243423
$ws.Cells.Item(2, 14).Value = 17.19057864341051
$ws.Cells.Item(2, 15).Value = 19.82791641170222
# Row 3
$ws.Cells.Item(3, 2).Value = 7.075783527932948
$ws.Cells.Item(3, 4).Value = 10.25299263611986
$ws.Cells.Item(3, 5).Value = 14.01658374581659
$ws.Cells.Item(3, 6).Value = 27.14206722156417
$ws.Cells.Item(3, 7).Value = 26.10151387948638
$ws.Cells.Item(3, 8).Value = 13.1571917940496
$ws.Cells.Item(3, 9).Value = 17.85781426329295
$ws.Cells.Item(3, 10).Value = 9.834813238555711
$ws.Cells.Item(3, 13).Value = 17.49164126111245
$ws.Cells.Item(3, 14).Value = 17.12657724550398
$ws.Cells.Item(3, 15).Value = 19.75164677367438
# Row 4
$ws.Cells.Item(4, 2).Value = 6.946824128499326
$ws.Cells.Item(4, 4).Value = 10.27229533786319
$ws.Cells.Item(4, 5).Value = 14.05690234757764
$ws.Cells.Item(4, 6).Value = 27.08144833209294
$ws.Cells.Item(4, 7).Value = 25.85320015632577
$ws.Cells.Item(4, 8).Value = 13.15951818302927
$ws.Cells.Item(4, 9).Value = 17.88140524998053
$ws.Cells.Item(4, 10).Value = 9.858895253336376
$ws.Cells.Item(4, 13).Value = 17.21164983886202
$ws.Cells.Item(4, 14).Value = 17.08914550076535
$ws.Cells.Item(4, 15).Value = 19.70983515849041
# Row 5
$ws.Cells.Item(5, 2).Value = 6.893673620465007
$ws.Cells.Item(5, 4).Value = 10.28054375151336
$ws.Cells.Item(5, 5).Value = 14.07394388543131
$ws.Cells.Item(5, 6).Value = 27.05861419314027
$ws.Cells.Item(5, 7).Value = 25.75362065279408
$ws.Cells.Item(5, 8).Value = 13.16112890367258
$ws.Cells.Item(5, 9).Value = 17.89195870376868
$ws.Cells.Item(5, 10).Value = 9.869019697965173
$ws.Cells.Item(5, 13).Value = 17.09642556692604
$ws.Cells.Item(5, 14).Value = 17.07437194821698
$ws.Cells.Item(5, 15).Value = 19.69407151781981
# Row 6
$ws.Cells.Item(6, 2).Value = 6.884814114210982
$ws.Cells.Item(6, 4).Value = 10.28193649067206
$ws.Cells.Item(6, 5).Value = 14.07681054926203
$ws.Cells.Item(6, 6).Value = 27.05493597573447
$ws.Cells.Item(6, 7).Value = 25.73718700300903
$ws.Cells.Item(6, 8).Value = 13.16143640136275
$ws.Cells.Item(6, 9).Value = 17.8937678966801
$ws.Cells.Item(6, 10).Value = 9.87071964568718
$ws.Cells.Item(6, 13).Value = 17.07722921115228
$ws.Cells.Item(6, 14).Value = 17.07194813008638
$ws.Cells.Item(6, 15).Value = 19.69153132806638
# Row 7
$ws.Cells.Item(7, 2).Value = 6.946109644810758
$ws.Cells.Item(7, 4).Value = 10.27240503050283
$ws.Cells.Item(7, 5).Value = 14.05712970010247
$ws.Cells.Item(7, 6).Value = 27.08113279345174
$ws.Cells.Item(7, 7).Value = 25.85185048218828
$ws.Cells.Item(7, 8).Value = 13.15953722189956
$ws.Cells.Item(7, 9).Value = 17.88154377038001
$ws.Cells.Item(7, 10).Value = 9.859030535789651
$ws.Cells.Item(7, 13).Value = 17.21010023784159
$ws.Cells.Item(7, 14).Value = 17.08894430082094
$ws.Cells.Item(7, 15).Value = 19.70961738642317
# Row 8
$ws.Cells.Item(8, 2).Value = 7.211141057027868
$ws.Cells.Item(8, 4).Value = 10.23370199296221
$ws.Cells.Item(8, 5).Value = 13.97564955386309
$ws.Cells.Item(8, 6).Value = 27.2130862379138
$ws.Cells.Item(8, 7).Value = 26.37158345622648
$ws.Cells.Item(8, 8).Value = 13.15697276477512
$ws.Cells.Item(8, 9).Value = 17.83585015283238
$ws.Cells.Item(8, 10).Value = 9.810176223420275
$ws.Cells.Item(8, 13).Value = 17.78614788552477
$ws.Cells.Item(8, 14).Value = 17.16813059625819
$ws.Cells.Item(8, 15).Value = 19.80058544243269
# Row 9
$ws.Cells.Item(9, 2).Value = 7.70670135305056
$ws.Cells.Item(9, 4).Value = 10.16987172916344
$ws.Cells.Item(9, 5).Value = 13.83515056613237
$ws.Cells.Item(9, 6).Value = 27.52921239607999
$ws.Cells.Item(9, 7).Value = 27.43049946469981
$ws.Cells.Item(9, 8).Value = 13.17279143768526
$ws.Cells.Item(9, 9).Value = 17.77580848747751
$ws.Cells.Item(9, 10).Value = 9.724133373030419
$ws.Cells.Item(9, 13).Value = 18.86942902924899
$ws.Cells.Item(9, 14).Value = 17.33771083450677
$ws.Cells.Item(9, 15).Value = 20.0181401809721
# Row 10
$ws.Cells.Item(10, 2).Value = 8.052458994604498
$ws.Cells.Item(10, 4).Value = 10.13033110705491
$ws.Cells.Item(10, 5).Value = 13.74367459360978
$ws.Cells.Item(10, 6).Value = 27.79496149997234
$ws.Cells.Item(10, 7).Value = 28.22394852766971
$ws.Cells.Item(10, 8).Value = 13.19704821453833
$ws.Cells.Item(10, 9).Value = 17.7496286730761
$ws.Cells.Item(10, 10).Value = 9.666827950204194
$ws.Cells.Item(10, 13).Value = 19.62950500750967
$ws.Cells.Item(10, 14).Value = 17.47031863386394
$ws.Cells.Item(10, 15).Value = 20.20087504577739
# Row 11
$ws.Cells.Item(11, 2).Value = 8.205116826577108
$ws.Cells.Item(11, 4).Value = 10.11394006120905
$ws.Cells.Item(11, 5).Value = 13.70460998202254
$ws.Cells.Item(11, 6).Value = 27.92276553024792
$ws.Cells.Item(11, 7).Value = 28.58643584849528
$ws.Cells.Item(11, 8).Value = 13.21080391964055
$ws.Cells.Item(11, 9).Value = 17.74158901307198
$ws.Cells.Item(11, 10).Value = 9.642033209488323
$ws.Cells.Item(11, 13).Value = 19.96609639682711
$ws.Cells.Item(11, 14).Value = 17.53223060446653
$ws.Cells.Item(11, 15).Value = 20.28873742602424
# Row 12
$ws.Cells.Item(12, 2).Value = 8.262210867226759
$ws.Cells.Item(12, 4).Value = 10.10796272587556
$ws.Cells.Item(12, 5).Value = 13.69018373435809
$ws.Cells.Item(12, 6).Value = 27.97212179882391
$ws.Cells.Item(12, 7).Value = 28.72377391269915
$ws.Cells.Item(12, 8).Value = 13.21640159507226
$ws.Cells.Item(12, 9).Value = 17.73909865082274
$ws.Cells.Item(12, 10).Value = 9.63282665674846
$ws.Cells.Item(12, 13).Value = 20.09213101566255
$ws.Cells.Item(12, 14).Value = 17.55588977726533
$ws.Cells.Item(12, 15).Value = 20.32266717306349
# Row 13
$ws.Cells.Item(13, 2).Value = 8.249947127958912
$ws.Cells.Item(13, 4).Value = 10.10923983954807
$ws.Cells.Item(13, 5).Value = 13.69327437296524
$ws.Cells.Item(13, 6).Value = 27.96144993144987
$ws.Cells.Item(13, 7).Value = 28.69419483928437
$ws.Cells.Item(13, 8).Value = 13.21517879384928
$ws.Cells.Item(13, 9).Value = 17.73961039232907
$ws.Cells.Item(13, 10).Value = 9.634801335674867
$ws.Cells.Item(13, 13).Value = 20.06505219027315
$ws.Cells.Item(13, 14).Value = 17.55078503714487
$ws.Cells.Item(13, 15).Value = 20.31533089835471
# Row 14
$ws.Cells.Item(14, 2).Value = 8.209828561336948
$ws.Cells.Item(14, 4).Value = 10.11344370038203
$ws.Cells.Item(14, 5).Value = 13.70341577581745
$ws.Cells.Item(14, 6).Value = 27.92680707148083
$ws.Cells.Item(14, 7).Value = 28.59773413330966
$ws.Cells.Item(14, 8).Value = 13.21125667009596
$ws.Cells.Item(14, 9).Value = 17.74137303929637
$ws.Cells.Item(14, 10).Value = 9.64127212249725
$ws.Cells.Item(14, 13).Value = 19.97649448160686
$ws.Cells.Item(14, 14).Value = 17.53417284138042
$ws.Cells.Item(14, 15).Value = 20.29151579085347
# Row 15
$ws.Cells.Item(15, 2).Value = 8.185160418079004
$ws.Cells.Item(15, 4).Value = 10.11604858887038
$ws.Cells.Item(15, 5).Value = 13.709675437033
$ws.Cells.Item(15, 6).Value = 27.90571123910094
$ws.Cells.Item(15, 7).Value = 28.53865416886853
$ws.Cells.Item(15, 8).Value = 13.20890479116949
$ws.Cells.Item(15, 9).Value = 17.74252479474729
$ws.Cells.Item(15, 10).Value = 9.645259441885637
$ws.Cells.Item(15, 13).Value = 19.9220617183733
$ws.Cells.Item(15, 14).Value = 17.52402490318635
$ws.Cells.Item(15, 15).Value = 20.27701336602491
# Row 16
$ws.Cells.Item(16, 2).Value = 8.042385039365517
$ws.Cells.Item(16, 4).Value = 10.13143442579335
$ws.Cells.Item(16, 5).Value = 13.74627885490994
$ws.Cells.Item(16, 6).Value = 27.78674541784935
$ws.Cells.Item(16, 7).Value = 28.20027726380094
$ws.Cells.Item(16, 8).Value = 13.1962037885533
$ws.Cells.Item(16, 9).Value = 17.7502317543711
$ws.Cells.Item(16, 10).Value = 9.668473934361007
$ws.Cells.Item(16, 13).Value = 19.60731411295707
$ws.Cells.Item(16, 14).Value = 17.466303313309
$ws.Cells.Item(16, 15).Value = 20.19522640278989
# Row 17
$ws.Cells.Item(17, 2).Value = 7.953575283238407
$ws.Cells.Item(17, 4).Value = 10.14128201049286
$ws.Cells.Item(17, 5).Value = 13.76938669624489
$ws.Cells.Item(17, 6).Value = 27.71551051641831
$ws.Cells.Item(17, 7).Value = 27.99298119499252
$ws.Cells.Item(17, 8).Value = 13.18910734627083
$ws.Cells.Item(17, 9).Value = 17.75594923110957
$ws.Cells.Item(17, 10).Value = 9.683041167038539
$ws.Cells.Item(17, 13).Value = 19.4117990463343
$ws.Cells.Item(17, 14).Value = 17.43128990677394
$ws.Cells.Item(17, 15).Value = 20.14624979731405
# Row 18
$ws.Cells.Item(18, 2).Value = 7.902061134329987
$ws.Cells.Item(18, 4).Value = 10.1470963172997
$ws.Cells.Item(18, 5).Value = 13.78291755438441
$ws.Cells.Item(18, 6).Value = 27.67519100266243
$ws.Cells.Item(18, 7).Value = 27.87390396533912
$ws.Cells.Item(18, 8).Value = 13.18528196913245
$ws.Cells.Item(18, 9).Value = 17.75960216158593
$ws.Cells.Item(18, 10).Value = 9.691539777142967
$ws.Cells.Item(18, 13).Value = 19.29848650206982
$ws.Cells.Item(18, 14).Value = 17.4113014420242
$ws.Cells.Item(18, 15).Value = 20.11852691192634
# Row 19
$ws.Cells.Item(19, 2).Value = 7.884546497166014
$ws.Cells.Item(19, 4).Value = 10.14909074271127
$ws.Cells.Item(19, 5).Value = 13.78754005637077
$ws.Cells.Item(19, 6).Value = 27.6616526663836
$ws.Cells.Item(19, 7).Value = 27.83361755085004
$ws.Cells.Item(19, 8).Value = 13.18403085729464
$ws.Cells.Item(19, 9).Value = 17.76090163343348
$ws.Cells.Item(19, 10).Value = 9.694437872573685
$ws.Cells.Item(19, 13).Value = 19.25997694974996
$ws.Cells.Item(19, 14).Value = 17.40455992845975
$ws.Cells.Item(19, 15).Value = 20.10921789291845
# Row 20
$ws.Cells.Item(20, 2).Value = 7.963074427171013
$ws.Cells.Item(20, 4).Value = 10.14021816884944
$ws.Cells.Item(20, 5).Value = 13.76690200304328
$ws.Cells.Item(20, 6).Value = 27.72302626246744
$ws.Cells.Item(20, 7).Value = 28.01503337904417
$ws.Cells.Item(20, 8).Value = 13.18983626545859
$ws.Cells.Item(20, 9).Value = 17.75530289944651
$ws.Cells.Item(20, 10).Value = 9.681478052245023
$ws.Cells.Item(20, 13).Value = 19.43270150220478
$ws.Cells.Item(20, 14).Value = 17.43500168596357
$ws.Cells.Item(20, 15).Value = 20.15141732183127
# Row 21
$ws.Cells.Item(21, 2).Value = 8.221632101798637
$ws.Cells.Item(21, 4).Value = 10.11220269282034
$ws.Cells.Item(21, 5).Value = 13.70042704654475
$ws.Cells.Item(21, 6).Value = 27.93695675861967
$ws.Cells.Item(21, 7).Value = 28.62606619912038
$ws.Cells.Item(21, 8).Value = 13.21239816659819
$ws.Cells.Item(21, 9).Value = 17.74084028964491
$ws.Cells.Item(21, 10).Value = 9.639366541595834
$ws.Cells.Item(21, 13).Value = 20.00254551047473
$ws.Cells.Item(21, 14).Value = 17.5390465377571
$ws.Cells.Item(21, 15).Value = 20.29849319027139
# Row 22
$ws.Cells.Item(22, 2).Value = 8.386432630881249
$ws.Cells.Item(22, 4).Value = 10.09523125850891
$ws.Cells.Item(22, 5).Value = 13.65911920099368
$ws.Cells.Item(22, 6).Value = 28.08234902756685
$ws.Cells.Item(22, 7).Value = 29.02576496948683
$ws.Cells.Item(22, 8).Value = 13.22940778640056
$ws.Cells.Item(22, 9).Value = 17.73461647906183
$ws.Cells.Item(22, 10).Value = 9.612908799372629
$ws.Cells.Item(22, 13).Value = 20.36662121237774
$ws.Cells.Item(22, 14).Value = 17.60828924862987
$ws.Cells.Item(22, 15).Value = 20.39844068234623
# Row 23
$ws.Cells.Item(23, 2).Value = 8.298872917151696
$ws.Cells.Item(23, 4).Value = 10.10416675735699
$ws.Cells.Item(23, 5).Value = 13.68097031940024
$ws.Cells.Item(23, 6).Value = 28.00425207495669
$ws.Cells.Item(23, 7).Value = 28.81245415768114
$ws.Cells.Item(23, 8).Value = 13.22012322007135
$ws.Cells.Item(23, 9).Value = 17.73764369815336
$ws.Cells.Item(23, 10).Value = 9.626932545865635
$ws.Cells.Item(23, 13).Value = 20.1731040124788
$ws.Cells.Item(23, 14).Value = 17.57122399436593
$ws.Cells.Item(23, 15).Value = 20.34475468251003
# Row 24
$ws.Cells.Item(24, 2).Value = 7.958781281541754
$ws.Cells.Item(24, 4).Value = 10.14069865571398
$ws.Cells.Item(24, 5).Value = 13.76802456704583
$ws.Cells.Item(24, 6).Value = 27.71962641498958
$ws.Cells.Item(24, 7).Value = 28.00506326346388
$ws.Cells.Item(24, 8).Value = 13.18950592823689
$ws.Cells.Item(24, 9).Value = 17.7555939661624
$ws.Cells.Item(24, 10).Value = 9.682184351039851
$ws.Cells.Item(24, 13).Value = 19.42325432292902
$ws.Cells.Item(24, 14).Value = 17.43332314968954
$ws.Cells.Item(24, 15).Value = 20.14907972905904
# Row 25
$ws.Cells.Item(25, 2).Value = 7.575618328287377
$ws.Cells.Item(25, 4).Value = 10.1858477527694
$ws.Cells.Item(25, 5).Value = 13.87109568737666
$ws.Cells.Item(25, 6).Value = 27.43769375697336
$ws.Cells.Item(25, 7).Value = 27.1406887807382
$ws.Cells.Item(25, 8).Value = 13.16628650446294
$ws.Cells.Item(25, 9).Value = 17.78889475117704
$ws.Cells.Item(25, 10).Value = 9.746369500046564
$ws.Cells.Item(25, 13).Value = 18.58215106567086
$ws.Cells.Item(25, 14).Value = 17.29037502529173
$ws.Cells.Item(25, 15).Value = 19.95518773915965

Write-Output "Updated loading_percent values for 380 kV case"
